$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# 1. Fix the "(in per cent)" -> "(in percent)" typo in C2
$ws.Range("C2").Value = "(in percent)"

# 2. Update the 2022 renewable-share figure (S5): 30 -> 29.9
$ws.Range("S5").Value = 29.9

# 3. Add the new 2023 column (T)
$ws.Range("S4").Copy()
$ws.Range("T4").PasteSpecial(-4122)
$ws.Range("T4").Value = 2023

$ws.Range("S5").Copy()
$ws.Range("T5").PasteSpecial(-4122)
$ws.Range("T5").Value = 29.5

$ws.Range("S6").Copy()
$ws.Range("T6").PasteSpecial(-4122)
$ws.Range("T6").Value = 12030.6

# 4. Match the new column-width formatting: columns D:T now share width 9
$ws.Range("D1:T1").EntireColumn.ColumnWidth = 9
